$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price + 1h volume change) scraped on
# Sat Jan 27 09:43:54 UTC 2024, plus the Toncoin/Cosmos row swap.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.623.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.264.10'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.21%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.99'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.47%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '91.81'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.33%  '
$ws.Range('E7').Value = '  +2.70%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.479'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.14'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +3.85%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.89'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0797'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.57'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.42%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.607.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.02%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.21'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.253.41'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.763'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +3.26%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '41.550.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.46'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.53%  '
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.91'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.67'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '240.30'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('E25').Value = '  +3.76%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +5.29%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '23.94'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.63%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.07'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '159.81'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '34.21'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +6.01%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.999'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.19'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +4.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0742'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.56%  '
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.105'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.82%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '16.60'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.07%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.80'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.94'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.65%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.047.35'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '19.52'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.09%  '
$ws.Range('E46').Value = '  +2.41%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.87'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.08%  '
$ws.Range('E48').Value = '  +6.96%  '
$ws.Range('E49').Value = '  +4.43%  '
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.72%  '
